{"js": "// Append \" a Tropy\" (as its own run) to the \"Jazykov\u00e9 prost\u0159edky\" heading,\n// turning it into \"Jazykov\u00e9 prost\u0159edky a Tropy\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst target = paragraphs.items.find((p) => p.text === \"Jazykov\u00e9 prost\u0159edky\");\n\nif (target) {\n  // Build a minimal Flat-OPC package so the inserted text lands in its own\n  // <w:r> (matching how the source document grew this heading), instead of\n  // being folded into the existing run's text.\n  const flatOpcRun =\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<?mso-application progid=\"Word.Document\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p><w:r>\" +\n    '<w:t xml:space=\"preserve\"> a Tropy</w:t>' +\n    \"</w:r></w:p></w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\";\n\n  const endRange = target.getRange(\"End\");\n  endRange.insertOoxml(flatOpcRun, \"End\");\n  await context.sync();\n}\n", "ps1": "# Append \" a Tropy\" (as its own run) to the \"Jazykov\u00e9 prost\u0159edky\" heading,\n# turning it into \"Jazykov\u00e9 prost\u0159edky a Tropy\".\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $full = $p.Range.Text\n    $text = $full.TrimEnd([char]13, [char]7)\n    if ($text -eq \"Jazykov\u00e9 prost\u0159edky\") {\n        # Range covering just the paragraph's text, excluding the trailing\n        # paragraph mark, so replacing it keeps <w:pPr> (the Heading 2 style)\n        # untouched.\n        $textRange = $d.Range($p.Range.Start, $p.Range.End - 1)\n\n        $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n            '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n            '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n            '<pkg:xmlData>' +\n            '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n            '<w:body><w:p>' +\n            '<w:r><w:t>Jazykov\u00e9 prost\u0159edky</w:t></w:r>' +\n            '<w:r><w:t xml:space=\"preserve\"> a Tropy</w:t></w:r>' +\n            '</w:p></w:body>' +\n            '</w:document>' +\n            '</pkg:xmlData></pkg:part></pkg:package>'\n\n        $textRange.InsertXML($xml)\n        break\n    }\n}\n"}
